# Weekly update: insert two new price records (week of 2022-11-11, serial 44876)
# for "Apio" at "Vega Monumental Concepción" above the existing row 270, pushing
# all subsequent rows (270-364) down by two (to 272-366).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 270, shifting existing data down to 272:366.
$ws.Rows("270:271").Insert()

# --- New row 270: Apio, Primera ---
$ws.Cells.Item(270, 1).Value = 11
$ws.Cells.Item(270, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(270, 3).Value = "Bíobío"
$ws.Cells.Item(270, 4).Value = 44876
$ws.Cells.Item(270, 5).Value = 8
$ws.Cells.Item(270, 6).Value = 100112017
$ws.Cells.Item(270, 7).Value = "Apio"
$ws.Cells.Item(270, 8).Value = "Americana (o)"
$ws.Cells.Item(270, 9).Value = "Primera"
$ws.Cells.Item(270, 10).Value = 350
$ws.Cells.Item(270, 11).Value = 7000
$ws.Cells.Item(270, 12).Value = 7500
$ws.Cells.Item(270, 13).Value = 7214
$ws.Cells.Item(270, 14).Value = "$/docena de matas"
$ws.Cells.Item(270, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(270, 16).Value = 1202
$ws.Cells.Item(270, 17).Value = 6
$ws.Cells.Item(270, 18).Value = "Hortaliza"

# --- New row 271: Apio, Segunda ---
$ws.Cells.Item(271, 1).Value = 11
$ws.Cells.Item(271, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(271, 3).Value = "Bíobío"
$ws.Cells.Item(271, 4).Value = 44876
$ws.Cells.Item(271, 5).Value = 8
$ws.Cells.Item(271, 6).Value = 100112017
$ws.Cells.Item(271, 7).Value = "Apio"
$ws.Cells.Item(271, 8).Value = "Americana (o)"
$ws.Cells.Item(271, 9).Value = "Segunda"
$ws.Cells.Item(271, 10).Value = 200
$ws.Cells.Item(271, 11).Value = 6000
$ws.Cells.Item(271, 12).Value = 6000
$ws.Cells.Item(271, 13).Value = 6000
$ws.Cells.Item(271, 14).Value = "$/docena de matas"
$ws.Cells.Item(271, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(271, 16).Value = 1000
$ws.Cells.Item(271, 17).Value = 6
$ws.Cells.Item(271, 18).Value = "Hortaliza"
